# Rename the drawing objects (wp:docPr / pic:cNvPr "name" attribute) for
# the three logo pictures living in the document's first-page footer,
# default footer and first-page header:
#   - PearsonLogo picture in the "first page" footer   : image2.png -> image1.png
#   - PearsonLogo picture in the default footer         : image2.png -> image1.png
#   - BTec_Logo-Orange picture in the "first page" header: image1.jpg -> image2.jpg
#
# InlineShape has no writable .Name property (matches real Word's object
# model), so each picture is briefly promoted to a floating Shape (which
# does expose .Name), renamed, then converted straight back to an inline
# shape so the drawing stays wp:inline exactly as before.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($headerFooter, $newName) {
    $inline = $headerFooter.Range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape()
}

# First-page footer (footer1.xml, docPr id="3") - Pearson logo.
$footerFirst = $sec.Footers.Item(2)
Rename-InlinePicture $footerFirst "image1.png"

# Default footer (footer2.xml, docPr id="2") - Pearson logo.
$footerDefault = $sec.Footers.Item(1)
Rename-InlinePicture $footerDefault "image1.png"

# First-page header (header1.xml, docPr id="1") - BTEC logo.
$headerFirst = $sec.Headers.Item(2)
Rename-InlinePicture $headerFirst "image2.jpg"
